$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New word-level data replacing the old phrase-level rows (2-5).
# Columns: A=text, B=x, C=y, D=width, E=height
$data = @(
    @("segregation ", 528, 644, 106, 23),
    @("Battle ", 673, 694, 55, 23),
    @("of ", 728, 694, 22, 23),
    @("Normandy. ", 750, 694, 98, 23),
    @("historically ", 527, 719, 95, 23),
    @("black ", 622, 719, 51, 23),
    @("college ", 673, 719, 67, 23),
    @("NAACP's ", 487, 993.1999999999999, 85, 23),
    @("boycotts ", 398, 1018.2, 77, 23),
    @("school ", 571, 1084.6, 61, 23),
    @("integration, ", 632, 1084.6, 101, 23),
    @("civil ", 269, 1176, 38, 23),
    @("rights ", 307, 1176, 53, 23),
    @("activists ", 360, 1176, 74, 23)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $item = $data[$i]
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
}
